# Update coin price/volume data and fix a row ordering swap (BOLO / CoinbaseStockToken)
# per the Jan 16 2023 08:44 UTC GitHub Actions symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "300.40"
$ws.Range("E2").Value = "1.46%"

# Row 3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.14%"

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.132"
$ws.Range("E4").Value = "0.38%"

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07958"
$ws.Range("E5").Value = "6.56%"

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "2.389"
$ws.Range("E6").Value = "41.85%"

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.68%"

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "3.845"
$ws.Range("E8").Value = "1.68%"

# Row 9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9114"
$ws.Range("E9").Value = "-1.93%"

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1722"
$ws.Range("E10").Value = "1.91%"

# Row 11
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07322"
$ws.Range("E11").Value = "1.89%"

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08050"
$ws.Range("E12").Value = "1.56%"

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03025"
$ws.Range("E13").Value = "0.74%"

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.41%"

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001532"
$ws.Range("E15").Value = "2.70%"

# Row 16
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006080"
$ws.Range("E16").Value = "-3.79%"

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.496"
$ws.Range("E17").Value = "1.23%"

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "2.241"
$ws.Range("E18").Value = "0.65%"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.85%"

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1354"
$ws.Range("E20").Value = "0.30%"

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "4.635"
$ws.Range("E21").Value = "1.08%"

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04622"
$ws.Range("E23").Value = "-0.78%"

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001264"
$ws.Range("E24").Value = "3.79%"

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004459"
$ws.Range("E25").Value = "0.80%"

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001194"
$ws.Range("E26").Value = "-8.21%"

# Row 27
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003442"
$ws.Range("E27").Value = "83.34%"

# Row 39
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01838"
$ws.Range("E39").Value = "10.83%"

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04526"
$ws.Range("E40").Value = "2.20%"

# Row 41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007301"
$ws.Range("E41").Value = "3.83%"

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1343"
$ws.Range("E42").Value = "1.23%"

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002173"
$ws.Range("E43").Value = "4.88%"

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.98%"

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006302"
$ws.Range("E45").Value = "5.31%"

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000753"
$ws.Range("E46").Value = "0.37%"

# Row 47
$ws.Range("B47:E47").NumberFormat = "@"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "0.006668"
$ws.Range("E47").Value = "-39.44%"

# Row 48
$ws.Range("B48:E48").NumberFormat = "@"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.8206"
$ws.Range("E48").Value = "15.31%"

# Row 49
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002108"
$ws.Range("E49").Value = "0.37%"

# Row 50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002008"
$ws.Range("E50").Value = "0.44%"
